$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.334.98"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.148.56"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.82"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.34"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.148.33"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("E10").Value = "  -4.70%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  -4.77%  "
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  -4.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.665.47"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.387.92"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.149.09"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.56"
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.71"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.85"
$ws.Range("E24").Value = "  -5.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.26"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("E29").Value = "  -4.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  -6.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.00"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -5.72%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.79"
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("E35").Value = "  -5.08%  "
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.51"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  +5.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0750"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "454.11"
$ws.Range("E40").Value = "  -9.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0401"
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.124"
$ws.Range("E42").Value = "  -5.02%  "
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.894.30"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("E45").Value = "  -7.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.32"
$ws.Range("E46").Value = "  -4.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.79"
$ws.Range("E47").Value = "  -4.94%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.55"
$ws.Range("E51").Value = "  -1.13%  "
